# Updated Test Suite table: renamed test case identifiers and added new test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Row 5/6: rename existing identifiers and swap their order
#   old B5 = TestCaseData_DataValidationWithURI, old B6 = TestCaseData_NavigationWithURI
#   new B5 = TestCaseData_DataValidationTestsWithURI, new B6 = TestCaseData_NavigationTestsWithURI
$ws.Range("B6").Value = "TestCaseData_NavigationTestsWithURI"
$ws.Range("B5").Value = "TestCaseData_DataValidationTestsWithURI"

# New rows 7-9
$ws.Range("A7").Value = "No"
$ws.Range("B7").Value = "TestCaseData_DataValidationTestsWithoutURI"
$ws.Range("A8").Value = "No"
$ws.Range("B8").Value = "TestCaseData_NavigationTestsWithoutURI"
$ws.Range("A9").Value = "No"
$ws.Range("B9").Value = "TestCaseData_CalculatorTestsWithURI"

# Resize the table and related ranges to include the new rows
$lo.Resize($ws.Range("A1:B9"))

# Widen column B to fit the longer identifier strings
$ws.Columns.Item(2).ColumnWidth = 45.67

# Move selection to the new last cell
[void]$ws.Range("A9").Select()
